$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column F header ("OSMO_DEF") is added after the existing OSMO_DESC
# column. Copy E1's formatting (bold font + border) onto F1, then set its
# text so it matches the look of the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "OSMO_DEF"
